$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Actual Results" column (G): top-file test bench results mirror the
# "Expected results / Outputs" column (F) for each of the 4 test cases.
$ws.Range("G2").Value2 = $ws.Range("F2").Value2
$ws.Range("G3").Value2 = $ws.Range("F3").Value2
$ws.Range("G4").Value2 = $ws.Range("F4").Value2
$ws.Range("G5").Value2 = $ws.Range("F5").Value2

# Copy the formatting of each "Expected results" cell onto its "Actual
# results" neighbour so the two columns line up visually.
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("F5").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- "Pass/Fail" column (H): every test case passed.
$ws.Range("H2:H5").Value2 = "PASS"
$ws.Range("H2:H5").VerticalAlignment = -4108

# --- "Git Commit" column (I): commit that introduced the top-file test bench.
$ws.Range("I2:I5").Value2 = "6c32a9f"
$ws.Range("I2:I5").VerticalAlignment = -4108
$ws.Range("I2:I5").HorizontalAlignment = -4131

# Widen the "Actual Results" column now that it holds real content.
$ws.Columns("G").ColumnWidth = 40.92

# Scroll the sheet over and move the selection, matching where the author
# left off after adding the new columns.
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("K4").Select()
